$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33.54101966249685
$ws.Range("B3").Value = 31.6227766016838

$ws.Range("A4:B4").Delete()
